$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Phase 1: cell values, entered in the same block order the original author
# used (LCD block, then Potentiometer block, then Button block) so that new
# shared-string entries land in the expected index order.
# ---------------------------------------------------------------------------

# LCD 16x2 I2C block
$ws.Range("G1").Value() = "LCD 16x2 I2C"
$ws.Range("G2").Value() = "LCD_pwr"
$ws.Range("G5").Value() = "LCD_gnd"
$ws.Range("D10").Value() = "A4"
$ws.Range("G10").Value() = "LCD_SDA"
$ws.Range("D11").Value() = "A5"
$ws.Range("G11").Value() = "LCD_SCL"

# Potentiometer block
$ws.Range("H1").Value() = "Potentiometer"
$ws.Range("D12").Value() = "A0"
$ws.Range("H2").Value() = "Pot_5v"
$ws.Range("H5").Value() = "Pot_GND"
$ws.Range("H12").Value() = "Pot_output"

# Button block
$ws.Range("I1").Value() = "Button"
$ws.Range("D13").Value() = "D7"
$ws.Range("I5").Value() = "Btn_GND"
$ws.Range("I13").Value() = "Btn_input"

# Two more rows appended below the previous bottom (12V/VIN/VIN and GND/GND/GND)
$ws.Range("C14").Value() = "12V"
$ws.Range("D14").Value() = "VIN"
$ws.Range("E14").Value() = "VIN"
$ws.Range("C15").Value() = "GND"
$ws.Range("D15").Value() = "GND"
$ws.Range("E15").Value() = "GND"

# ---------------------------------------------------------------------------
# Phase 2: formatting
# ---------------------------------------------------------------------------

# G2:I15 block -> thin border all round, horizontal-center only
$gi = $ws.Range("G2:I15")
$gi.Borders.LineStyle = 1
$gi.HorizontalAlignment = -4108

# C2:F15 block -> thin border all round, horizontal+vertical center
$cf = $ws.Range("C2:F15")
$cf.Borders.LineStyle = 1
$cf.HorizontalAlignment = -4108
$cf.VerticalAlignment = -4108

# A few G/H/I cells match the C:F look (center+center) rather than the
# plain G:I look (center only)
foreach ($addr in @("H2","H5","I5")) {
    $c = $ws.Range($addr)
    $c.VerticalAlignment = -4108
}

# Header row C1:F1 -> bold, highlight fill, border left+right+top thin (no
# bottom), center+center
foreach ($addr in @("C1","D1","E1","F1")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Interior.Color = 15652797
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(10).LineStyle = 1
    $c.Borders.Item(8).LineStyle = 1
    $c.Borders.Item(9).LineStyle = -4142
}

# Header row G1:I1 -> bold, highlight fill, border right-thin only,
# center+center
foreach ($addr in @("G1","H1","I1")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Interior.Color = 15652797
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.Borders.Item(10).LineStyle = 1
}

# ---------------------------------------------------------------------------
# Phase 3: new column widths for G, H, I
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 13.0
$ws.Columns.Item(9).ColumnWidth = 10.166666666666666

# ---------------------------------------------------------------------------
# Phase 4: selection, matching the saved cursor position
# ---------------------------------------------------------------------------
$ws.Range("L11").Select()
